# Washington_sign_images.xlsx edit: added and deleted some images
#
# Summary of the change (derived from the OOXML diff):
#  - Two rows were removed from the middle of the table:
#      * the duplicate "A-2829033" row that pointed at image 0033.png (E=0)
#      * the "A-2829014" row that pointed at image 0034.png (E=2)
#    (these accident-id/image rows are dropped entirely)
#  - The remaining "A-2829033" row's data_frame_row_number (col B) was
#    corrected from 26 to 16.
#  - The "Number" column (col A) is a simple 1..N row counter, so after the
#    deletion the later rows were renumbered to stay contiguous.
#  - Five new rows were appended at the end of the table for a new accident
#    "A-2828757" (data_frame_row_number 20) with images 0053-0057.png.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the two obsolete rows (sheet rows 29 and 30) ---------------
# Row 29: Number=28, B=26, C=A-2829033, D=0033.png, E=0
# Row 30: Number=29, B=25, C=A-2829014, D=0034.png, E=2
$ws.Range("A29:E30").EntireRow.Delete()

# --- 2. Fix the data_frame_row_number for the remaining A-2829033 row -----
# (now at sheet row 28, the former row 27, untouched by the delete above)
$ws.Range("B28").Value = 16

# --- 3. Renumber the "Number" column so it stays a contiguous 1..N --------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# --- 4. Append the five new rows for accident A-2828757 -------------------
$newRows = @(
    @(20, "A-2828757", "0053.png", 33),
    @(20, "A-2828757", "0054.png", 34),
    @(20, "A-2828757", "0055.png", 37),
    @(20, "A-2828757", "0056.png", 33),
    @(20, "A-2828757", "0057.png", 37)
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp, recompute after renumbering
foreach ($row in $newRows) {
    $lastRow = $lastRow + 1
    $ws.Cells.Item($lastRow, 1).Value = $lastRow - 1
    $ws.Cells.Item($lastRow, 2).Value = $row[0]
    $ws.Cells.Item($lastRow, 3).Value = $row[1]
    $ws.Cells.Item($lastRow, 4).Value = $row[2]
    $ws.Cells.Item($lastRow, 5).Value = $row[3]
}

# --- 5. Restore the view state (scroll position / selection) --------------
$ws.Activate()
$ws.Range("B47").Select()
